# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (columns E and F, rows 16-79)
# was re-sorted from newest-period-first to oldest-period-first
# (chronological ascending order instead of descending). Reversing the
# 64-row block for both columns reproduces the new period labels and
# their corresponding arrears amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$endRow = 79
$n = $endRow - $startRow + 1

# Snapshot the current (descending) values first, since we will
# overwrite the range in place.
$eVals = @()
$fVals = @()
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $eVals += $ws.Cells.Item($r, 5).Value2
    $fVals += $ws.Cells.Item($r, 6).Value2
}

# Write them back in reverse (ascending chronological) order.
for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $srcIdx = $n - 1 - $i
    $ws.Cells.Item($r, 5).Value = $eVals[$srcIdx]
    $ws.Cells.Item($r, 6).Value = $fVals[$srcIdx]
}
